$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Add new log row 27 (mirrors the structure of existing rows)
$ws.Range("A27").Value = "Wanneer zijn jullie open?"
$ws.Range("B27").Value = "mailmind.test@zohomail.eu"
$ws.Range("C27").Value = "Testmail #1: Wanneer zijn jullie open?"
$ws.Range("D27").Value = "Openingstijden / Locatie"
$ws.Range("E27").Value = "Beste klant,`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$ws.Range("F27").Value = "2025-06-26 22:35:27"
$ws.Range("G27").Value = "Ja"
$ws.Range("H27").Value = "Nee"
$ws.Range("I27").Value = "Ja"

# Reset auto-calculated row height back to standard (matches non-wrapped rows)
$ws.Rows.Item(27).EntireRow.AutoFit()

# Extend conditional formatting ranges to cover the new row
$dRule = $ws.Range("D2:D26").FormatConditions.Item(1)
$dRule.ModifyAppliesToRange($ws.Range("D2:D27"))

$gRule = $ws.Range("G2:G26").FormatConditions.Item(1)
$gRule.ModifyAppliesToRange($ws.Range("G2:G27"))

$hRule = $ws.Range("H2:H26").FormatConditions.Item(1)
$hRule.ModifyAppliesToRange($ws.Range("H2:H27"))

$iRule = $ws.Range("I2:I26").FormatConditions.Item(1)
$iRule.ModifyAppliesToRange($ws.Range("I2:I27"))

# Update Dashboard count for "Openingstijden / Locatie" category (3 -> 4)
$dashboard.Range("B3").Value = 4
